$p = $ppt.ActivePresentation

$oldStyleId = "{7FF2914C-5757-49FB-8CC2-641F3AF1693F}"
$newStyleId = "{B111C9E1-FE86-41AE-A90A-032E88E1AC1B}"

for ($slideIdx = 1; $slideIdx -le $p.Slides.Count; $slideIdx++) {
    $slide = $p.Slides.Item($slideIdx)
    for ($shapeIdx = 1; $shapeIdx -le $slide.Shapes.Count; $shapeIdx++) {
        $shape = $slide.Shapes.Item($shapeIdx)
        if ($shape.HasTable) {
            $tbl = $shape.Table
            if ($tbl.Style -eq $oldStyleId) {
                $tbl.ApplyStyle($newStyleId)
            }
        }
    }
}

# Re-theme the deck from the custom "Integral / Red Violet" design to the
# standard Office theme palette (dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink).
$s1 = $p.Slides.Item(1)
$themeColors = $s1.ThemeColorScheme
$themeColors.Item(1).RGB  = 0         # dk1      000000
$themeColors.Item(2).RGB  = 16777215  # lt1      FFFFFF
$themeColors.Item(3).RGB  = 6968388   # dk2      44546A
$themeColors.Item(4).RGB  = 15132391  # lt2      E7E6E6
$themeColors.Item(5).RGB  = 13998939  # accent1  5B9BD5
$themeColors.Item(6).RGB  = 3243501   # accent2  ED7D31
$themeColors.Item(7).RGB  = 10855845  # accent3  A5A5A5
$themeColors.Item(8).RGB  = 49407     # accent4  FFC000
$themeColors.Item(9).RGB  = 12874308  # accent5  4472C4
$themeColors.Item(10).RGB = 4697456   # accent6  70AD47
$themeColors.Item(11).RGB = 12673797  # hlink    0563C1
$themeColors.Item(12).RGB = 7491477   # folHlink 954F72
